# Insert a new weekly price record as row 12 (pushing the existing rows
# 12-48 down to 13-49), matching the pattern already used throughout this
# "Camote" (Vega Modelo de Temuco) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12:48 down one, creating a blank (but format-inherited) row 12.
$ws.Rows("12:12").Insert()

# Populate the new row 12 with the new record.
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value = "La Araucanía"
$ws.Cells.Item(12, 4).Value = 44613
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 100114002
$ws.Cells.Item(12, 7).Value = "Camote"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 11).Value = 18000
$ws.Cells.Item(12, 12).Value = 18000
$ws.Cells.Item(12, 13).Value = 18000
$ws.Cells.Item(12, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 900
$ws.Cells.Item(12, 17).Value = 20
$ws.Cells.Item(12, 18).Value = "Hortaliza"
